$d = $word.ActiveDocument

# The document currently ends with an empty list-item paragraph
# (style CustomStyle, numId 1, ilvl 2) that was left blank. Fill it
# in with the first new bullet text, then append the remaining
# bullets as new paragraphs after it, each with its own list level.
# Note: Word's ListLevelNumber is 1-based (ilvl + 1).
$p = $d.Paragraphs.Last
$p.Range.ListFormat.ListLevelNumber = 3
$p.Range.Text = "Each provided stat for players"

$null = $p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.ListLevelNumber = 3
$p.Range.Text = "Which team they have been drafted to/free agent"

$null = $p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.ListLevelNumber = 3
$p.Range.Text = "Order in which they have been drafted"

$null = $p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.ListLevelNumber = 2
$p.Range.Text = "Team class will track each team’s roster"

$null = $p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.ListLevelNumber = 3
$p.Range.Text = "Players filling each slot"

$null = $p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.ListLevelNumber = 3
$p.Range.Text = "Methods to draft or identify draftees by other participants"

$null = $p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.ListLevelNumber = 2
$p.Range.Text = "Methods class will include:"

$null = $p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.ListLevelNumber = 3
$p.Range.Text = "Other methods to manage the databases and other fantasy draft functions"

